$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("G9","G10","G11","G13","G14","G19","H19","G21","H21")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("C8").Value = 77
$ws.Range("C9").Value = 60
$ws.Range("G9").Value = "15360.00"
$ws.Range("C10").Value = 11
$ws.Range("G10").Value = "5192.00"
$ws.Range("C11").Value = 66
$ws.Range("G11").Value = "43692.00"
$ws.Range("C12").Value = 71
$ws.Range("C13").Value = 90
$ws.Range("G13").Value = "12240.00"
$ws.Range("C14").Value = 38
$ws.Range("G14").Value = "874.00"
$ws.Range("C15").Value = 71
$ws.Range("C16").Value = 55
$ws.Range("C17").Value = 37
$ws.Range("G19").Value = "77358.00"
$ws.Range("H19").Value = "77358.00"
$ws.Range("G21").Value = "77358.00"
$ws.Range("H21").Value = "77358.00"
